$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column B width to match column A (target stored width 15.42578125).
# NB: Excel COM's ColumnWidth is points-based (internally stored as
# pixels = chars * 12, then re-quantized to 1/6-pt increments on readback),
# so it can only land on multiples of 1/6. 14.65 is the input that rounds to
# the closest achievable stored width (15.5) to the target 15.42578125.
$ws.Columns.Item(2).ColumnWidth = 14.65

$ws.Cells.Item(1, 1).Value = -0.34640275735485204
$ws.Cells.Item(1, 2).Value = 0.3456671421520383
$ws.Cells.Item(2, 1).Value = -0.2550774472426731
$ws.Cells.Item(2, 2).Value = 0.25277667161775774
$ws.Cells.Item(3, 1).Value = -0.21908224146656607
$ws.Cells.Item(3, 2).Value = 0.2177970787250736
$ws.Cells.Item(4, 1).Value = -0.2057970789153849
$ws.Cells.Item(4, 2).Value = 0.20461525485141507
$ws.Cells.Item(5, 1).Value = -0.1986152555374332
$ws.Cells.Item(5, 2).Value = 0.19622292254282758
$ws.Cells.Item(6, 1).Value = -0.0987725748874162
$ws.Cells.Item(6, 2).Value = 0.0986210092927311
$ws.Cells.Item(7, 1).Value = -0.07862101013088107
$ws.Cells.Item(7, 2).Value = 0.07825758510678149
$ws.Cells.Item(8, 1).Value = -0.05825758595398689
$ws.Cells.Item(8, 2).Value = 0.05796714621216381
$ws.Cells.Item(9, 1).Value = -0.051967146949039034
$ws.Cells.Item(9, 2).Value = 0.05172783274432735
$ws.Cells.Item(10, 1).Value = -0.04572783349060927
$ws.Cells.Item(10, 2).Value = 0.04569632062029427
$ws.Cells.Item(11, 1).Value = -0.04119632135393658
$ws.Cells.Item(11, 2).Value = 0.04114514633729627
$ws.Cells.Item(12, 1).Value = -0.035145147086030004
$ws.Cells.Item(12, 2).Value = 0.03499096012776137
$ws.Cells.Item(13, 1).Value = -0.028990960885580286
$ws.Cells.Item(13, 2).Value = 0.028953177830272736
$ws.Cells.Item(14, 1).Value = -0.0169531786424848
$ws.Cells.Item(14, 2).Value = 0.01694123485302068
$ws.Cells.Item(15, 1).Value = -0.010941235615106848
$ws.Cells.Item(15, 2).Value = 0.010935063887658636
$ws.Cells.Item(16, 1).Value = -0.004935064650964716
$ws.Cells.Item(16, 2).Value = 0.004932402657806145
$ws.Cells.Item(17, 1).Value = 0.0010675965779389784
$ws.Cells.Item(17, 2).Value = -0.0010674999390012374
$ws.Cells.Item(18, 1).Value = -0.08058969558740259
$ws.Cells.Item(18, 2).Value = 0.08052783770920868
$ws.Cells.Item(19, 1).Value = -0.07152783837926702
$ws.Cells.Item(19, 2).Value = 0.07104140326520625
$ws.Cells.Item(20, 1).Value = -0.018013298244566656
$ws.Cells.Item(20, 2).Value = 0.018004258310373444
$ws.Cells.Item(21, 1).Value = -0.009004258997572734
$ws.Cells.Item(21, 2).Value = 0.008999999312166551
$ws.Cells.Item(22, 1).Value = -0.0246944308401158
$ws.Cells.Item(22, 2).Value = 0.02467707723047674
$ws.Cells.Item(23, 1).Value = -0.015677077920950744
$ws.Cells.Item(23, 2).Value = 0.01566299079040867
$ws.Cells.Item(24, 1).Value = -0.04212464602765831
$ws.Cells.Item(24, 2).Value = 0.04199999901793383
$ws.Cells.Item(25, 1).Value = -0.09145034836970822
$ws.Cells.Item(25, 2).Value = 0.09133648151319207
$ws.Cells.Item(26, 1).Value = -0.0853364822299092
$ws.Cells.Item(26, 2).Value = 0.08519503869520406
$ws.Cells.Item(27, 1).Value = -0.07831458418206116
$ws.Cells.Item(27, 2).Value = 0.07784619758069367
$ws.Cells.Item(28, 1).Value = -0.07184619831526806
$ws.Cells.Item(28, 2).Value = 0.07154461620314656
$ws.Cells.Item(29, 1).Value = -0.05954461699932878
$ws.Cells.Item(29, 2).Value = 0.059418196120880395
$ws.Cells.Item(30, 1).Value = -0.042169966388336366
$ws.Cells.Item(30, 2).Value = 0.04201907678625272
$ws.Cells.Item(31, 1).Value = -0.027019077620252574
$ws.Cells.Item(31, 2).Value = 0.027000600596004176
$ws.Cells.Item(32, 1).Value = -0.006000601483020063
$ws.Cells.Item(32, 2).Value = 0.005999999242169096
